$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Sheet is protected; unprotect to edit then restore protection afterward
$ws.Unprotect()

# Update the confidential disclaimer text date from 2021-03-22 to 2021-03-23
$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-23 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-6
$ws.Range("D2").Value = 0.2490016186540238
$ws.Range("E2").Value = -0.02084405558414826

$ws.Range("D3").Value = 0.2456656150064469
$ws.Range("E3").Value = -0.01397561700862326

$ws.Range("D4").Value = 0.2505057047397971
$ws.Range("E4").Value = -0.006218245241525699

$ws.Range("D5").Value = 0.2548270615997321
$ws.Range("E5").Value = -0.007377157159794479

$ws.Range("D6").Value = 0.9999999999999999
$ws.Range("E6").Value = -0.01206113731574854

# Restore sheet protection
$ws.Protect()
